$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp string (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 18 de Octubre de 2020 a las 17:52"

# --- Refresh case numbers for the countries whose stats changed ---
# Columns: A=Pais B=Casos totales C=Nuevos casos D=Casos activos E=Recuperados F=Casos criticos G=Muertes hoy H=Muertes

# Estados Unidos (row 4)
$ws.Range("B4").Value = 8349951
$ws.Range("C4").Value = 7286
$ws.Range("D4").Value = 5437959
$ws.Range("E4").Value = 2687621
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 89
$ws.Range("H4").Value = 224371

# India (row 5)
$ws.Range("B5").Value = 7531825
$ws.Range("C5").Value = 39098
$ws.Range("D5").Value = 6636342
$ws.Range("E5").Value = 781048
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 371
$ws.Range("H5").Value = 114435

# Reino Unido (row 14)
$ws.Range("B14").Value = 722409
$ws.Range("C14").Value = 16982
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 67
$ws.Range("H14").Value = 43646

# Italia (row 19)
$ws.Range("B19").Value = 414241
$ws.Range("C19").Value = 11705
$ws.Range("D19").Value = 251461
$ws.Range("E19").Value = 126237
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 69
$ws.Range("H19").Value = 36543

# Canada (row 31)
$ws.Range("B31").Value = 198075
$ws.Range("C31").Value = 1754
$ws.Range("D31").Value = 167089
$ws.Range("E31").Value = 21229
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 11
$ws.Range("H31").Value = 9757

# Republica Dominicana (row 41)
$ws.Range("B41").Value = 121347
$ws.Range("C41").Value = 422
$ws.Range("D41").Value = 98207
$ws.Range("E41").Value = 20941
$ws.Range("F41").Value = 0
$ws.Range("G41").Value = 4
$ws.Range("H41").Value = 2199

# Guatemala (row 48)
$ws.Range("B48").Value = 101360
$ws.Range("C48").Value = 332
$ws.Range("D48").Value = 90610
$ws.Range("E48").Value = 7220
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 15
$ws.Range("H48").Value = 3530

# Japon (row 51)
$ws.Range("B51").Value = 92656
$ws.Range("C51").Value = 593
$ws.Range("D51").Value = 85485
$ws.Range("E51").Value = 5501
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 9
$ws.Range("H51").Value = 1670

# Moldavia (row 59)
$ws.Range("B59").Value = 67050
$ws.Range("C59").Value = 398
$ws.Range("D59").Value = 47842
$ws.Range("E59").Value = 17624
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 15
$ws.Range("H59").Value = 1584

# Singapur (row 65) - only Casos activos / Recuperados change
$ws.Range("D65").Value = 57807
$ws.Range("E65").Value = 76

# Kenia (row 75)
$ws.Range("B75").Value = 44881
$ws.Range("C75").Value = 685
$ws.Range("D75").Value = 31857
$ws.Range("E75").Value = 12192
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 7
$ws.Range("H75").Value = 832

# Jordania moves ahead of Serbia (rows 78/79 swap + Jordania gets new figures)
$ws.Range("A78").Value = "Jordania"
$ws.Range("B78").Value = 37573
$ws.Range("C78").Value = 1520
$ws.Range("D78").Value = 6912
$ws.Range("E78").Value = 30316
$ws.Range("F78").Value = 0
$ws.Range("G78").Value = 15
$ws.Range("H78").Value = 345

$ws.Range("A79").Value = "Serbia"
$ws.Range("B79").Value = 36160
$ws.Range("C79").Value = 214
$ws.Range("D79").Value = 31536
$ws.Range("E79").Value = 3848
$ws.Range("F79").Value = 0
$ws.Range("G79").Value = 2
$ws.Range("H79").Value = 776

# Grecia moves ahead of Corea del Sur (rows 88/89 swap + Grecia gets new figures)
$ws.Range("A88").Value = "Grecia"
$ws.Range("B88").Value = 25370
$ws.Range("C88").Value = 438
$ws.Range("D88").Value = 9989
$ws.Range("E88").Value = 14872
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 9
$ws.Range("H88").Value = 509

$ws.Range("A89").Value = "Corea del Sur"
$ws.Range("B89").Value = 25199
$ws.Range("C89").Value = 91
$ws.Range("D89").Value = 23312
$ws.Range("E89").Value = 1443
$ws.Range("F89").Value = 0
$ws.Range("G89").Value = 1
$ws.Range("H89").Value = 444

# Zambia (row 98) - only first four data columns change
$ws.Range("B98").Value = 15853
$ws.Range("C98").Value = 64
$ws.Range("D98").Value = 15005
$ws.Range("E98").Value = 502

# Jamaica (row 115)
$ws.Range("B115").Value = 8274
$ws.Range("C115").Value = 79
$ws.Range("D115").Value = 3859
$ws.Range("E115").Value = 4244
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 3
$ws.Range("H115").Value = 171

# Cuba (row 122) - only first four data columns change
$ws.Range("B122").Value = 6220
$ws.Range("C122").Value = 50
$ws.Range("D122").Value = 5768
$ws.Range("E122").Value = 327

# Republica de Yibuti (row 127) - only first four data columns change
$ws.Range("B127").Value = 5459
$ws.Range("C127").Value = 7
$ws.Range("D127").Value = 5375
$ws.Range("E127").Value = 23

# Sudan del Sur (row 152) - Casos totales, Nuevos casos, Recuperados change
$ws.Range("B152").Value = 2842
$ws.Range("C152").Value = 25
$ws.Range("E152").Value = 1497

# Burundi (row 177) - Casos totales, Nuevos casos, Recuperados change
$ws.Range("B177").Value = 542
$ws.Range("C177").Value = 6
$ws.Range("E177").Value = 44
